$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Stash the four distinct cell-formats used across the table (label
# cell, text cell, numeric-qty cell, plain bordered/empty cell) into a
# scratch area off to the side so we can restore them after Merge()
# (which rewrites borders on the ranges it touches).
# ------------------------------------------------------------------
$ws.Range("A3").Copy()
$ws.Range("Z1").PasteSpecial(-4122)   # style "3" - label cell (A column)
$ws.Range("B3").Copy()
$ws.Range("Z2").PasteSpecial(-4122)   # style "4" - text cell (B/D column)
$ws.Range("C3").Copy()
$ws.Range("Z3").PasteSpecial(-4122)   # style "5" - numeric qty cell (C column)
$ws.Range("E3").Copy()
$ws.Range("Z4").PasteSpecial(-4122)   # style "6" - plain bordered/empty cell

# ------------------------------------------------------------------
# Row 1 - title cell gets text
# ------------------------------------------------------------------
$ws.Range("A1").Value = "dupa"

# ------------------------------------------------------------------
# Cells whose text happens to look like a plain number still need to
# stay genuinely text-typed (matching the source inline strings), so
# force a text number-format on them right before assigning - Excel
# would otherwise silently re-interpret "33" etc. as the number 33.
# ------------------------------------------------------------------
$textCells = @("B3", "B5", "B6", "B7", "D7", "D8")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# ------------------------------------------------------------------
# Row 3 - edit existing values
# ------------------------------------------------------------------
$ws.Range("A3").Value = "M3"
$ws.Range("B3").Value = "33"
$ws.Range("C3").Value = 12

# ------------------------------------------------------------------
# Row 4 (new)
# ------------------------------------------------------------------
$ws.Range("C4").Value = 21
$ws.Range("D4").Value = "bialy"

# ------------------------------------------------------------------
# Row 5 (new)
# ------------------------------------------------------------------
$ws.Range("A5").Value = "D9"
$ws.Range("B5").Value = "12"
$ws.Range("C5").Value = 12
$ws.Range("D5").Value = "styropian"

# ------------------------------------------------------------------
# Row 6 (new)
# ------------------------------------------------------------------
$ws.Range("A6").Value = "Statyw drewniany"
$ws.Range("B6").Value = "34"
$ws.Range("C6").Value = 34
$ws.Range("D6").Value = "biały"

# ------------------------------------------------------------------
# Row 7 (new)
# ------------------------------------------------------------------
$ws.Range("A7").Value = "Statyw metalowy"
$ws.Range("B7").Value = "46"
$ws.Range("C7").Value = 12
$ws.Range("D7").Value = "90"

# ------------------------------------------------------------------
# Row 8 (new)
# ------------------------------------------------------------------
$ws.Range("C8").Value = 34
$ws.Range("D8").Value = "45"

# ------------------------------------------------------------------
# Merges
# ------------------------------------------------------------------
$ws.Range("A3:A4").Merge()
$ws.Range("B3:B4").Merge()
$ws.Range("A5").Merge()
$ws.Range("B5").Merge()
$ws.Range("A6").Merge()
$ws.Range("B6").Merge()
$ws.Range("A7:A8").Merge()
$ws.Range("B7:B8").Merge()

# ------------------------------------------------------------------
# Restore proper cell formats (Merge() rewrites borders on touched
# cells, so reapply the canonical per-column look afterwards).
# ------------------------------------------------------------------
$ws.Range("Z1").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A7").PasteSpecial(-4122)

$ws.Range("Z2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D8").PasteSpecial(-4122)

$ws.Range("Z3").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("C8").PasteSpecial(-4122)

$ws.Range("Z4").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("E7").PasteSpecial(-4122)
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("E8").PasteSpecial(-4122)

# ------------------------------------------------------------------
# Clean up scratch cells
# ------------------------------------------------------------------
$ws.Range("Z1:Z4").Clear()
